$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4886

$ws.Range("H80").Value = 570.2857
$ws.Range("J80").Value = 549
$ws.Range("L80").Value = 1647
$ws.Range("N80").Value = -3643

$ws.Range("H83").Value = 570.2857
$ws.Range("J83").Value = 549
$ws.Range("L83").Value = 4941
$ws.Range("N83").Value = -14925

$ws.Range("H103").Value = 2262.4285
$ws.Range("J103").Value = 1061.5
$ws.Range("L103").Value = 3184.5
$ws.Range("N103").Value = -4356.5

$ws.Range("H106").Value = 5523
$ws.Range("I106").Value = 5523
$ws.Range("K106").Value = 5523
$ws.Range("M106").Value = -4892

$ws.Range("H107").Value = 1177.3636
$ws.Range("I107").Value = 1414.1111
$ws.Range("K107").Value = 1414.1111
$ws.Range("M107").Value = 505.8888999999999

$ws.Range("H138").Value = 3329.5454
$ws.Range("I138").Value = 1638.2
$ws.Range("J138").Value = 4739
$ws.Range("K138").Value = 4914.6
$ws.Range("L138").Value = 14217
$ws.Range("M138").Value = 225.3999999999996
$ws.Range("N138").Value = -24497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 112.6
$ws.Range("I5").Value = 115.75
$ws.Range("K5").Value = 115.75
$ws.Range("M5").Value = -3.75

$ws.Range("H24").Value = 5005177.5
$ws.Range("J24").Value = 5005177.5
$ws.Range("L24").Value = 5005177.5
$ws.Range("N24").Value = -5005925.5

$ws.Range("H32").Value = 569
$ws.Range("I32").Value = 576.5599999999999
$ws.Range("K32").Value = 576.5599999999999
$ws.Range("M32").Value = -289.5599999999999

$ws.Range("H45").Value = 2210.7778
$ws.Range("I45").Value = 1343.3334
$ws.Range("K45").Value = 1343.3334
$ws.Range("M45").Value = -966.3334

$ws.Range("H63").Value = 5821.4443
$ws.Range("I63").Value = 1770.7142
$ws.Range("K63").Value = 1770.7142
$ws.Range("M63").Value = -1084.7142

$ws.Range("H66").Value = 5821.4443
$ws.Range("I66").Value = 1770.7142
$ws.Range("K66").Value = 8853.571
$ws.Range("M66").Value = -5421.571

$ws.Range("H97").Value = 2330.8333
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("N97").Value = -2992

$ws.Range("H100").Value = 5005177.5
$ws.Range("J100").Value = 5005177.5
$ws.Range("L100").Value = 5005177.5
$ws.Range("N100").Value = -5007341.5

$ws.Range("H110").Value = 432.63635
$ws.Range("I110").Value = 666.3333
$ws.Range("J110").Value = 152.2
$ws.Range("K110").Value = 666.3333
$ws.Range("L110").Value = 152.2
$ws.Range("M110").Value = 1378.6667
$ws.Range("N110").Value = -4242.2

$ws.Range("H122").Value = 3280.8572
$ws.Range("I122").Value = 3280.8572
$ws.Range("K122").Value = 9842.571599999999
$ws.Range("M122").Value = -7392.571599999999

$ws.Range("H132").Value = 1773.75
$ws.Range("I132").Value = 1771.3636
$ws.Range("K132").Value = 5314.0908
$ws.Range("M132").Value = -2784.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 112.6
$ws.Range("I4").Value = 115.75
$ws.Range("K4").Value = 115.75
$ws.Range("M4").Value = -0.75

$ws.Range("H22").Value = 5380
$ws.Range("I22").Value = 5380
$ws.Range("K22").Value = 5380
$ws.Range("M22").Value = -5207

$ws.Range("H86").Value = 3028
$ws.Range("I86").Value = 1608.2941
$ws.Range("K86").Value = 1608.2941
$ws.Range("M86").Value = -485.2941000000001

$ws.Range("H89").Value = 3028
$ws.Range("I89").Value = 1608.2941
$ws.Range("K89").Value = 8041.4705
$ws.Range("M89").Value = -2425.4705

$ws.Range("H107").Value = 6185.5
$ws.Range("I107").Value = 5812.5
$ws.Range("K107").Value = 5812.5
$ws.Range("M107").Value = -3892.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9100.799999999999
$ws.Range("I62").Value = 9000
$ws.Range("K62").Value = 9000
$ws.Range("M62").Value = -8376

$ws.Range("H65").Value = 9100.799999999999
$ws.Range("I65").Value = 9000
$ws.Range("K65").Value = 45000
$ws.Range("M65").Value = -41880

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H105").Value = 1070.3334
$ws.Range("I105").Value = 1070.3334
$ws.Range("K105").Value = 1070.3334
$ws.Range("M105").Value = 676.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1000
$ws.Range("J52").Value = 1000
$ws.Range("L52").Value = 3000
$ws.Range("N52").Value = -3532

$ws.Range("H114").Value = 485
$ws.Range("I114").Value = 197
$ws.Range("K114").Value = 591
$ws.Range("M114").Value = 2663

$ws.Range("H121").Value = 1338.1666
$ws.Range("I121").Value = 676.3333
$ws.Range("K121").Value = 2028.9999
$ws.Range("M121").Value = -718.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2869.2727
$ws.Range("I80").Value = 3286.25
$ws.Range("K80").Value = 3286.25
$ws.Range("M80").Value = -2288.25

$ws.Range("H83").Value = 2869.2727
$ws.Range("I83").Value = 3286.25
$ws.Range("K83").Value = 16431.25
$ws.Range("M83").Value = -11439.25

$ws.Range("H102").Value = 2708.2727
$ws.Range("I102").Value = 1897.6666
$ws.Range("K102").Value = 1897.6666
$ws.Range("M102").Value = -275.6666

$ws.Range("H113").Value = 4719.8335
$ws.Range("I113").Value = 3904.8333
$ws.Range("K113").Value = 3904.8333
$ws.Range("M113").Value = -1734.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1217.1428
$ws.Range("I55").Value = 1110
$ws.Range("K55").Value = 1110
$ws.Range("M55").Value = -937

$ws.Range("H68").Value = 4246.5
$ws.Range("I68").Value = 2328.6667
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 2328.6667
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -1579.6667
$ws.Range("N68").Value = -11498

$ws.Range("H71").Value = 4246.5
$ws.Range("I71").Value = 2328.6667
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 11643.3335
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -7899.333500000001
$ws.Range("N71").Value = -57488

$ws.Range("H82").Value = 3891.25
$ws.Range("I82").Value = 2569.111
$ws.Range("J82").Value = 5591.143
$ws.Range("K82").Value = 2569.111
$ws.Range("L82").Value = 5591.143
$ws.Range("M82").Value = -2208.111
$ws.Range("N82").Value = -6313.143

$ws.Range("H85").Value = 3891.25
$ws.Range("I85").Value = 2569.111
$ws.Range("J85").Value = 5591.143
$ws.Range("K85").Value = 2569.111
$ws.Range("L85").Value = 5591.143
$ws.Range("M85").Value = -1321.111
$ws.Range("N85").Value = -8087.143

$ws.Range("H100").Value = 4554.222
$ws.Range("I100").Value = 1831.3334
$ws.Range("K100").Value = 1831.3334
$ws.Range("M100").Value = -1290.3334

$ws.Range("H122").Value = 2997.6
$ws.Range("I122").Value = 2997.6
$ws.Range("K122").Value = 8992.799999999999
$ws.Range("M122").Value = -6542.799999999999

$ws.Range("H132").Value = 6421.8
$ws.Range("I132").Value = 5777.25
$ws.Range("K132").Value = 17331.75
$ws.Range("M132").Value = -14801.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1012.125
$ws.Range("I100").Value = 822.61536
$ws.Range("K100").Value = 1645.23072
$ws.Range("M100").Value = -1104.23072

$ws.Range("H122").Value = 4502.25
$ws.Range("I122").Value = 4502.25
$ws.Range("K122").Value = 13506.75
$ws.Range("M122").Value = -11056.75

$ws.Range("H136").Value = 3177.88
$ws.Range("I136").Value = 2074.8125
$ws.Range("K136").Value = 6224.4375
$ws.Range("M136").Value = -4224.4375
